{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the five existing installation-steps list paragraphs by their\n// current text so the script is resilient to exact indices.\nconst items = paragraphs.items;\n\nfunction findByStart(startText) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(startText) === 0) {\n      return items[i];\n    }\n  }\n  throw new Error(\"Paragraph starting with '\" + startText + \"' not found\");\n}\n\nconst pJar = findByStart(\"Copy and paste jar files into Niagara Modules folder\");\nconst pDrag = findByStart(\"Drag an instance of the Schedule Parser and the Schedule Writer\");\nconst pCsv = findByStart(\"Copy and paste CSV of schedules into\");\nconst pGentle = findByStart(\"Run the Gentle Write action on the Schedule Writer\");\n\n// 1. \"Copy and paste jar files...\" -> \"Create a new program object\"\npJar.insertText(\"Create a new program object\", \"Replace\");\n\n// 2. \"Drag an instance...\" -> \"Copy and paste the supplied code for the Schedule Parser into a program object\"\npDrag.insertText(\"Copy and paste the supplied code for the Schedule Parser into a program object\", \"Replace\");\n\n// 3. \"Copy and paste CSV of schedules into Schedule Parser.\" becomes \"Save and compile\"\n//    and is followed by two brand-new list paragraphs before the (modified) CSV line.\npCsv.insertText(\"Save and compile\", \"Replace\");\nawait context.sync();\n\nlet newPara1 = pCsv.insertParagraph(\"Repeat 1-3 for the Schedule Writer\", \"After\");\nnewPara1.style = \"List Paragraph\";\nawait context.sync();\nnewPara1.attachToList(4, 0);\nawait context.sync();\n\nlet newPara2 = newPara1.insertParagraph(\"Copy and paste CSV of schedules into schedule parser.\", \"After\");\nnewPara2.style = \"List Paragraph\";\nawait context.sync();\nnewPara2.attachToList(4, 0);\nawait context.sync();\n\n// 5. Remove the second sentence/run appended to the Gentle Write bullet.\nconst gentleRange = pGentle.getRange(\"Whole\");\ngentleRange.load(\"text\");\nawait context.sync();\n\nconst gentleKeep = \"Run the Gentle Write action on the Schedule Writer (This will take some time)\";\ngentleRange.insertText(gentleKeep, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByPrefix($doc, $prefix) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        if ($t.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Set-ParagraphText($doc, $index, $newText) {\n    $p = $doc.Paragraphs.Item($index)\n    $r = $p.Range\n    $full = $doc.Range($r.Start, $r.End)\n    $full.Text = $newText\n}\n\n# --- Locate the five existing installation-step bullets ------------------\n$idxJar = Get-ParagraphIndexByPrefix $d \"Copy and paste jar files into Niagara Modules folder\"\n$idxDrag = Get-ParagraphIndexByPrefix $d \"Drag an instance of the Schedule Parser and the Schedule Writer\"\n$idxCsv = Get-ParagraphIndexByPrefix $d \"Copy and paste CSV of schedules into\"\n$idxGentle = Get-ParagraphIndexByPrefix $d \"Run the Gentle Write action on the Schedule Writer\"\n\n# 1. \"Copy and paste jar files...\" -> \"Create a new program object\"\nSet-ParagraphText $d $idxJar \"Create a new program object\"\n\n# 2. \"Drag an instance...\" -> \"Copy and paste the supplied code for the Schedule Parser into a program object\"\nSet-ParagraphText $d $idxDrag \"Copy and paste the supplied code for the Schedule Parser into a program object\"\n\n# 3. \"Copy and paste CSV of schedules into Schedule Parser.\" -> \"Save and compile\"\nSet-ParagraphText $d $idxCsv \"Save and compile\"\n\n# 4/5. Insert two brand-new list bullets after the (now \"Save and compile\") paragraph,\n#      reusing the same list (numId) as the rest of the bulleted list.\n$pCsv = $d.Paragraphs.Item($idxCsv)\n$listTemplate = $pCsv.Range.ListFormat.ListTemplate\n\n$pCsv.Range.InsertParagraphAfter()\n$idxNew1 = $idxCsv + 1\n$pNew1 = $d.Paragraphs.Item($idxNew1)\n$pNew1.Range.Text = \"Repeat 1-3 for the Schedule Writer\"\n$pNew1.Style = \"List Paragraph\"\n$pNew1.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, 1, 1)\n\n$pNew1.Range.InsertParagraphAfter()\n$idxNew2 = $idxNew1 + 1\n$pNew2 = $d.Paragraphs.Item($idxNew2)\n$pNew2.Range.Text = \"Copy and paste CSV of schedules into schedule parser.\"\n$pNew2.Style = \"List Paragraph\"\n$pNew2.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, 1, 1)\n\n# 6. Remove the second sentence/run appended to the Gentle Write bullet\n#    (re-locate it, since paragraph indices shifted after the inserts above).\n$idxGentle2 = Get-ParagraphIndexByPrefix $d \"Run the Gentle Write action on the Schedule Writer\"\nSet-ParagraphText $d $idxGentle2 \"Run the Gentle Write action on the Schedule Writer (This will take some time)\"\n"}
